$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 9: hours = 0, comments = new string
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = "Fue semana de examenes, y no tube mucho tiempo"

# Fill in week labels for rows 11-13
$ws.Range("B11").Value = "4 Octubre - 11 Octubre"
$ws.Range("B12").Value = "11 Octubre - 18- Octubre"
$ws.Range("B13").Value = "18 Octubre - 25 Octubre"

# Update selection to B13 as last active cell
$ws.Range("B13").Select()
